$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.1030763333333333
$ws.Range("H2").Value = 0.309229
$ws.Range("I2").Value = 0.01126512502660735
$ws.Range("J2").Value = 0.01126512502660735
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7004376666666666
$ws.Range("N2").Value = 2.101313
$ws.Range("O2").Value = 0.04511966030063898
$ws.Range("P2").Value = 0.04511966030063898
$ws.Range("Q2").Value = 0.07219854640855554
$ws.Range("R2").Value = 0.6497869176769999
$ws.Range("S2").Value = 0.0005082786144447504
$ws.Range("T2").Value = 0.0005082786144447504

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.1030763333333333
$ws.Range("H3").Value = 0.309229
$ws.Range("I3").Value = 0.01126512502660735
$ws.Range("J3").Value = 0.01126512502660735
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.237305
$ws.Range("N3").Value = 0.711915
$ws.Range("O3").Value = 0.01528632952964618
$ws.Range("P3").Value = 0.01528632952964618
$ws.Range("Q3").Value = 0.02446052928166666
$ws.Range("R3").Value = 0.220144763535
$ws.Range("S3").Value = 0.0001722024133493842
$ws.Range("T3").Value = 0.0001722024133493842

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.1030763333333333
$ws.Range("H4").Value = 0.309229
$ws.Range("I4").Value = 0.01126512502660735
$ws.Range("J4").Value = 0.01126512502660735
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.58625866666667
$ws.Range("N4").Value = 43.758776
$ws.Range("O4").Value = 0.9395940101697148
$ws.Range("P4").Value = 0.9395940101697148
$ws.Range("Q4").Value = 1.503498060411555
$ws.Range("R4").Value = 13.531482543704
$ws.Range("S4").Value = 0.01058464399881322
$ws.Range("T4").Value = 0.01058464399881322

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.046962666666667
$ws.Range("H5").Value = 27.140888
$ws.Range("I5").Value = 0.9887348749733926
$ws.Range("J5").Value = 0.9887348749733927
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7004376666666666
$ws.Range("N5").Value = 2.101313
$ws.Range("O5").Value = 0.04511966030063898
$ws.Range("P5").Value = 0.04511966030063898
$ws.Range("Q5").Value = 6.336833420660445
$ws.Range("R5").Value = 57.031500785944
$ws.Range("S5").Value = 0.04461138168619423
$ws.Range("T5").Value = 0.04461138168619423

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 9.046962666666667
$ws.Range("H6").Value = 27.140888
$ws.Range("I6").Value = 0.9887348749733926
$ws.Range("J6").Value = 0.9887348749733927
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.237305
$ws.Range("N6").Value = 0.711915
$ws.Range("O6").Value = 0.01528632952964618
$ws.Range("P6").Value = 0.01528632952964618
$ws.Range("Q6").Value = 2.146889475613333
$ws.Range("R6").Value = 19.32200528052
$ws.Range("S6").Value = 0.01511412711629679
$ws.Range("T6").Value = 0.0151141271162968

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 9.046962666666667
$ws.Range("H7").Value = 27.140888
$ws.Range("I7").Value = 0.9887348749733926
$ws.Range("J7").Value = 0.9887348749733927
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.58625866666667
$ws.Range("N7").Value = 43.758776
$ws.Range("O7").Value = 0.9395940101697148
$ws.Range("P7").Value = 0.9395940101697148
$ws.Range("Q7").Value = 131.9613376036764
$ws.Range("R7").Value = 1187.652038433088
$ws.Range("S7").Value = 0.9290093661709016
$ws.Range("T7").Value = 0.9290093661709017
